# Automatische test-sync: 2025-07-23 21:14:50
# Appends a new test-mail log entry (row 9) to the "Logs" sheet, extends the
# conditional-formatting ranges that cover the data rows, and bumps the
# matching "Aantal" counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append the new log row (row 9) -----------------------------------
$newRow = 9

$logs.Range("A$newRow").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B$newRow").Value = "mailmind.test@zohomail.eu"
$logs.Range("C$newRow").Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Range("D$newRow").Value = "Openingstijden / Locatie"
$logs.Range("E$newRow").Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F$newRow").Value = "2025-07-23 21:14:43"
$logs.Range("G$newRow").Value = "Ja"
$logs.Range("H$newRow").Value = "Nee"
$logs.Range("I$newRow").Value = "Ja"
$logs.Range("J$newRow").Value = "Ja"

# --- Extend conditional formatting ranges to include the new row -------
# Each sqref block (D, G, H, I, J) needs its AppliesTo range pushed from
# row 8 down to row 9; updating one rule in a block updates the whole
# shared sqref for that block.
$logs.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D9"))
$logs.Range("G2:G8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G9"))
$logs.Range("H2:H8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H9"))
$logs.Range("I2:I8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I9"))
$logs.Range("J2:J8").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J9"))

# --- Update the Dashboard summary count ---------------------------------
$dash.Range("B2").Value = 8
